# Clean up labels in raw data:
#  - Insert a header row ("kurzname" / "langname") at the top of the
#    "labels" sheet, shifting all existing rows down by one.
#  - Add a "Sonstige Anmerkungen" label next to the existing "Anmerkung"
#    row (now row 33, column B).
#  - Leave the "labels" sheet as the active sheet/tab (it was the last
#    sheet worked on), with B33 selected.

$wb = $excel.ActiveWorkbook
$labelsSheet = $wb.Worksheets.Item("labels")

# Insert new header row at top of labels sheet, pushing data down.
$labelsSheet.Rows.Item(1).Insert()
$labelsSheet.Range("A1").Value = "kurzname"
$labelsSheet.Range("B1").Value = "langname"

# Fill in the long name for the "Anmerkung" row (now row 33).
$labelsSheet.Range("B33").Value = "Sonstige Anmerkungen"

# Make "labels" the active sheet/tab, with B33 selected.
$labelsSheet.Activate()
$labelsSheet.Range("B33").Select()
